$d = $word.ActiveDocument

# Merge " as an actuated " + "pendulum" runs (with proofErr markers around
# "pendulum") into a single run reading " as an actuated pendulum", removing
# the proofErr gramStart/gramEnd markers in between.
$d.Content.Find.Execute("as an actuated pendulum", $true, $false, $false, $false, $false, $true, 1, $false, "as an actuated pendulum", 2) | Out-Null
